$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 132304.8
$ws.Range("I15").Value = 132304.8
$ws.Range("K15").Value = 396914.4
$ws.Range("M15").Value = -396745.4

$ws.Range("H132").Value = 229559.83
$ws.Range("I132").Value = 248871.94
$ws.Range("K132").Value = 746615.8200000001
$ws.Range("M132").Value = -744085.8200000001

$ws.Range("H137").Value = 17242344
$ws.Range("I137").Value = 20834046
$ws.Range("J137").Value = 2177.8
$ws.Range("K137").Value = 62502138
$ws.Range("L137").Value = 6533.400000000001
$ws.Range("M137").Value = -62499588
$ws.Range("N137").Value = -11633.4

$ws.Range("H138").Value = 8334840
$ws.Range("I138").Value = 948510.4399999999
$ws.Range("J138").Value = 333333340
$ws.Range("K138").Value = 2845531.32
$ws.Range("L138").Value = 1000000020
$ws.Range("M138").Value = -2840391.32
$ws.Range("N138").Value = -1000010300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 85499.336
$ws.Range("I2").Value = 145341.72
$ws.Range("J2").Value = 1720
$ws.Range("K2").Value = 145341.72
$ws.Range("L2").Value = 1720
$ws.Range("M2").Value = -145228.72
$ws.Range("N2").Value = -1946

$ws.Range("H32").Value = 18073.908
$ws.Range("I32").Value = 2686.2952
$ws.Range("K32").Value = 2686.2952
$ws.Range("M32").Value = -2399.2952

$ws.Range("H116").Value = 85499.336
$ws.Range("I116").Value = 145341.72
$ws.Range("J116").Value = 1720
$ws.Range("K116").Value = 145341.72
$ws.Range("L116").Value = 1720
$ws.Range("M116").Value = -143047.72
$ws.Range("N116").Value = -6308

$ws.Range("H122").Value = 2593.913
$ws.Range("I122").Value = 2624.1875
$ws.Range("J122").Value = 2524.7144
$ws.Range("K122").Value = 7872.5625
$ws.Range("L122").Value = 7574.1432
$ws.Range("M122").Value = -5422.5625
$ws.Range("N122").Value = -12474.1432

$ws.Range("H132").Value = 1639.8923
$ws.Range("I132").Value = 1205.4746
$ws.Range("J132").Value = 5911.6665
$ws.Range("K132").Value = 3616.4238
$ws.Range("L132").Value = 17734.9995
$ws.Range("M132").Value = -1086.4238
$ws.Range("N132").Value = -22794.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 85499.336
$ws.Range("I3").Value = 145341.72
$ws.Range("J3").Value = 1720
$ws.Range("K3").Value = 145341.72
$ws.Range("L3").Value = 1720
$ws.Range("M3").Value = -145227.72
$ws.Range("N3").Value = -1948

$ws.Range("H94").Value = 1117.2858
$ws.Range("I94").Value = 952.94446
$ws.Range("J94").Value = 2103.3333
$ws.Range("K94").Value = 952.94446
$ws.Range("L94").Value = 2103.3333
$ws.Range("M94").Value = -501.94446
$ws.Range("N94").Value = -3005.3333

$ws.Range("H99").Value = 2227.182
$ws.Range("I99").Value = 2149.9
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2149.9
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -651.9000000000001
$ws.Range("N99").Value = -5996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7813689
$ws.Range("I99").Value = 8929500
$ws.Range("J99").Value = 3014
$ws.Range("K99").Value = 8929500
$ws.Range("L99").Value = 3014
$ws.Range("M99").Value = -8928002
$ws.Range("N99").Value = -6010

$ws.Range("H126").Value = 7813689
$ws.Range("I126").Value = 8929500
$ws.Range("J126").Value = 3014
$ws.Range("K126").Value = 26788500
$ws.Range("L126").Value = 9042
$ws.Range("M126").Value = -26786030
$ws.Range("N126").Value = -13982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1057.0217
$ws.Range("I5").Value = 711.3714
$ws.Range("J5").Value = 2156.818
$ws.Range("K5").Value = 2134.1142
$ws.Range("L5").Value = 6470.454000000001
$ws.Range("M5").Value = -2022.1142
$ws.Range("N5").Value = -6694.454000000001

$ws.Range("H17").Value = 755.55554
$ws.Range("J17").Value = 755.55554
$ws.Range("L17").Value = 2266.66662
$ws.Range("N17").Value = -2604.66662

$ws.Range("H35").Value = 2500
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 4000
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = -2712
$ws.Range("N35").Value = -12576

$ws.Range("H107").Value = 860.96
$ws.Range("J107").Value = 404.8
$ws.Range("L107").Value = 1214.4
$ws.Range("N107").Value = -5054.4

$ws.Range("H113").Value = 14706789
$ws.Range("J113").Value = 25001090
$ws.Range("L113").Value = 75003270
$ws.Range("N113").Value = -75007610

$ws.Range("H122").Value = 492.625
$ws.Range("I122").Value = 268.33334
$ws.Range("J122").Value = 866.44446
$ws.Range("K122").Value = 2415.00006
$ws.Range("L122").Value = 7798.00014
$ws.Range("M122").Value = 34.9999399999997
$ws.Range("N122").Value = -12698.00014

$ws.Range("H135").Value = 1057.0217
$ws.Range("I135").Value = 711.3714
$ws.Range("J135").Value = 2156.818
$ws.Range("K135").Value = 6402.3426
$ws.Range("L135").Value = 19411.362
$ws.Range("M135").Value = -3867.3426
$ws.Range("N135").Value = -24481.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4800
$ws.Range("I12").Value = 4600
$ws.Range("K12").Value = 4600
$ws.Range("M12").Value = -4460

$ws.Range("H80").Value = 2626.3635
$ws.Range("I80").Value = 2598.3333
$ws.Range("J80").Value = 2660
$ws.Range("K80").Value = 2598.3333
$ws.Range("L80").Value = 2660
$ws.Range("M80").Value = -1600.3333
$ws.Range("N80").Value = -4656

$ws.Range("H83").Value = 2626.3635
$ws.Range("I83").Value = 2598.3333
$ws.Range("J83").Value = 2660
$ws.Range("K83").Value = 12991.6665
$ws.Range("L83").Value = 13300
$ws.Range("M83").Value = -7999.666499999999
$ws.Range("N83").Value = -23284

$ws.Range("H102").Value = 2720.875
$ws.Range("I102").Value = 2678
$ws.Range("J102").Value = 2849.5
$ws.Range("K102").Value = 2678
$ws.Range("L102").Value = 2849.5
$ws.Range("M102").Value = -1056
$ws.Range("N102").Value = -6093.5

$ws.Range("H122").Value = 1011354.75
$ws.Range("I122").Value = 1390051.5
$ws.Range("J122").Value = 1496.6666
$ws.Range("K122").Value = 4170154.5
$ws.Range("L122").Value = 4489.9998
$ws.Range("M122").Value = -4167704.5
$ws.Range("N122").Value = -9389.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2988.6333
$ws.Range("I40").Value = 1145.8182
$ws.Range("J40").Value = 4055.5264
$ws.Range("K40").Value = 1145.8182
$ws.Range("L40").Value = 4055.5264
$ws.Range("M40").Value = -1009.8182
$ws.Range("N40").Value = -4327.526400000001

$ws.Range("H55").Value = 332
$ws.Range("I55").Value = 216
$ws.Range("J55").Value = 399.66666
$ws.Range("K55").Value = 216
$ws.Range("L55").Value = 399.66666
$ws.Range("M55").Value = -43
$ws.Range("N55").Value = -745.66666

$ws.Range("H92").Value = 32000
$ws.Range("J92").Value = 32000
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992

$ws.Range("H93").Value = 1620.8
$ws.Range("I93").Value = 1400
$ws.Range("K93").Value = 1400
$ws.Range("M93").Value = -152

$ws.Range("H100").Value = 2843815.2
$ws.Range("I100").Value = 8930705
$ws.Range("J100").Value = 3266.6667
$ws.Range("K100").Value = 8930705
$ws.Range("L100").Value = 3266.6667
$ws.Range("M100").Value = -8930164
$ws.Range("N100").Value = -4348.6667

$ws.Range("H122").Value = 3642.8572
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 3769.2307
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 11307.6921
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -16207.6921

$ws.Range("H132").Value = 4690.3774
$ws.Range("I132").Value = 4794.282
$ws.Range("J132").Value = 4400.9287
$ws.Range("K132").Value = 14382.846
$ws.Range("L132").Value = 13202.7861
$ws.Range("M132").Value = -11852.846
$ws.Range("N132").Value = -18262.7861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 4884.5
$ws.Range("J31").Value = 4884.5
$ws.Range("L31").Value = 4884.5
$ws.Range("N31").Value = -5580.5

$ws.Range("H81").Value = 1820685.5
$ws.Range("I81").Value = 2502030
$ws.Range("K81").Value = 5004060
$ws.Range("M81").Value = -5002999

$ws.Range("H84").Value = 1820685.5
$ws.Range("I84").Value = 2502030
$ws.Range("K84").Value = 25020300
$ws.Range("M84").Value = -25014996

$ws.Range("H122").Value = 93454.91
$ws.Range("I122").Value = 202000.8
$ws.Range("K122").Value = 606002.3999999999
$ws.Range("M122").Value = -603552.3999999999

$ws.Range("H132").Value = 7694098.5
$ws.Range("I132").Value = 11112880
$ws.Range("J132").Value = 1840.45
$ws.Range("K132").Value = 33338640
$ws.Range("L132").Value = 5521.35
$ws.Range("M132").Value = -33336110
$ws.Range("N132").Value = -10581.35

$ws.Range("H136").Value = 13433.025
$ws.Range("I136").Value = 14771.112
$ws.Range("K136").Value = 44313.336
$ws.Range("M136").Value = -41763.336
